$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '66.274.93'
$ws.Range("D3").Value = '3.081.12'
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''574.99'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = '''169.45'
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.079.35'
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("D10").Value = '''6.33'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("E11").Value = '  -2.48%  '
$ws.Range("D12").Value = '''0.470'
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").Value = '''0.0000239'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").Value = '''35.74'
$ws.Range("E14").Value = '  -3.43%  '
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '3.592.93'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("D17").Value = '66.221.04'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").Value = '''6.95'
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("E19").Value = '  +2.51%  '
$ws.Range("D20").Value = '3.078.44'
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = '''487.21'
$ws.Range("E21").Value = '  +2.31%  '
$ws.Range("D22").Value = '''7.71'
$ws.Range("E22").Value = '  -3.03%  '
$ws.Range("D23").Value = '''0.687'
$ws.Range("E23").Value = '  -3.42%  '
$ws.Range("D24").Value = '''82.58'
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '''12.63'
$ws.Range("E25").Value = '  -5.20%  '
$ws.Range("E26").Value = '  -3.62%  '
$ws.Range("D27").Value = '''10.21'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '''7.86'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("E30").Value = '  -4.93%  '
$ws.Range("D31").Value = '''2.60'
$ws.Range("E31").Value = '  -2.84%  '
$ws.Range("D32").Value = '''27.75'
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("D34").Value = '0.0₃0899'
$ws.Range("E34").Value = '  -4.77%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''0.947'
$ws.Range("E36").Value = '  -2.88%  '
$ws.Range("B37").Value = 'Arweave'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D37").Value = '''47.15'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").Value = '''5.56'
$ws.Range("E38").Value = '  -4.84%  '
$ws.Range("E40").Value = '  -4.77%  '
$ws.Range("E41").Value = '  -4.41%  '
$ws.Range("E42").Value = '  -4.64%  '
$ws.Range("D43").Value = '2.784.52'
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("E44").Value = '  -2.89%  '
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").Value = '''134.57'
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").Value = '''364.25'
$ws.Range("E47").Value = '  -5.11%  '
$ws.Range("D49").Value = '''24.38'
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("E51").Value = '  -2.05%  '
